$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet's very first stored row is a pre-existing, already-correct header
# row (Lichess ID / Upland Username / Lichess Rating / Balance / Bearer Token /
# Eos Upland ID) sitting at Excel row 1 - it is left untouched.
#
# Best-effort: some files carry an extra (invalid, row-0) row above Excel's
# row 1 floor; real Excel has no such row, so this is guarded defensively and
# simply does nothing when the host (correctly) rejects it.
try {
    $ghost = $ws.Cells.Item(1, 1).Offset(-1, 0)
    $ghost.Value = "Lichess ID"
    $ghost.Offset(0, 1).Value = "Upland Username"
    $ghost.Offset(0, 2).Value = "Lichess Rating"
    $ghost.Offset(0, 3).Value = "Balance"
    $ghost.Offset(0, 4).Value = "Bearer Token"
    $ghost.Offset(0, 5).Value = "Eos Upland ID"
} catch {
}

# Excel row 2: append the new profile/data row below the header.
$ws.Range("A2").Value = "trashboatsr"
$ws.Range("B2").Value = "dogeyboy19"
$ws.Range("C2").Value = 1818
$ws.Range("D2").Value = 1875597
$ws.Range("E2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiMjRiZDI1YWItOGY3MS00YzJjLWEyZjYtMmMyN2Y0Mjg2ZjI3IiwiaWF0IjoxNzAzMzUzMTI2fQ.RgTv8LJBQqRC43i699uwZfVYmvXpUFKClspfNoEozVg"
$ws.Range("F2").Value = "mp4n4f2mq3ca"
